$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Column C ("Fitness") values were updated for Generations 0-92 (rows 2-94).
$ws.Range("C2:C49").Value = 7534
$ws.Range("C50:C60").Value = 7295
$ws.Range("C61:C94").Value = 7293
